$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row (rows 2-11).
# Bump the value by one day: 45181 (2023-09-12) -> 45182 (2023-09-13).
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45182
}
